# Fix length of Picoblade cable on solenoid driver BOM, and add the
# OC_RS column with the new part's RS Components order code.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the CABLE row (row 14): shorten cable length 100mm -> 50mm,
# and update the Molex part number.
$ws.Range("E14").Value = "Picoblade cable - 5 way 50mm"
$ws.Range("G14").Value = "15134-0500"

# New "OC_RS" header column (I1)
$ws.Range("I1").Value = "OC_RS"

# Add the RS order code for the cable.
$ws.Range("H14").Value = "125-0735"

# Match the saved selection state from the authored workbook.
$ws.Range("A14:H14").Select()
